$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy header style (bold, border, centered) from an existing header cell (AC1)
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill team record data for all data rows (2-55)
$ws.Range("AD2:AD55").Value = 76
$ws.Range("AE2:AE55").Value = 86
$ws.Range("AF2:AF55").Value = 0
